$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: stamp a cell with the same look (font/border/alignment) as the
# existing header-style cells (A1..H1 / A2..E2 use style index 1: bold,
# thin box border, centered, top-aligned) without touching its value.
# ---------------------------------------------------------------------------
function Copy-Format($srcCell, $dstCell) {
    $srcCell.Copy() | Out-Null
    $dstCell.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

# ---------------------------------------------------------------------------
# Row 2 — update the existing measurement row
# ---------------------------------------------------------------------------

# B2: 488 (number) -> "488" (text). Leading apostrophe forces text entry;
# re-stamp the formatting afterwards so it keeps the table's normal look.
$ws.Cells.Item(2, 2).Value = "'488"
Copy-Format $ws.Cells.Item(2, 1) $ws.Cells.Item(2, 2)

$ws.Cells.Item(2, 3).Value = 50

# D2: date-looking text -> force as text (not an Excel date) then restore style.
$ws.Cells.Item(2, 4).Value = "'2022-07-22"
Copy-Format $ws.Cells.Item(1, 4) $ws.Cells.Item(2, 4)

$ws.Cells.Item(2, 5).Value = "11:36"

$ws.Cells.Item(2, 7).Value = 50
$ws.Cells.Item(2, 8).Value = 30

# ---------------------------------------------------------------------------
# Row 3 — new row (A3/B3 stay empty but keep the header-style formatting)
# ---------------------------------------------------------------------------
$ws.Cells.Item(3, 3).Value = 100
$ws.Cells.Item(3, 4).Value = "'2022-07-22"
$ws.Cells.Item(3, 5).Value = "11:36"
$ws.Cells.Item(3, 6).Value = 0
$ws.Cells.Item(3, 7).Value = 100
$ws.Cells.Item(3, 8).Value = 30

Copy-Format $ws.Cells.Item(2, 1) $ws.Cells.Item(3, 1)
Copy-Format $ws.Cells.Item(2, 2) $ws.Cells.Item(3, 2)
Copy-Format $ws.Cells.Item(2, 3) $ws.Cells.Item(3, 3)
Copy-Format $ws.Cells.Item(1, 4) $ws.Cells.Item(3, 4)
Copy-Format $ws.Cells.Item(2, 5) $ws.Cells.Item(3, 5)

# ---------------------------------------------------------------------------
# Row 4 — new row
# ---------------------------------------------------------------------------
$ws.Cells.Item(4, 2).Value = "'561"
$ws.Cells.Item(4, 3).Value = 50
$ws.Cells.Item(4, 4).Value = "'2022-07-22"
$ws.Cells.Item(4, 5).Value = "11:36"
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 40
$ws.Cells.Item(4, 8).Value = 25

Copy-Format $ws.Cells.Item(2, 1) $ws.Cells.Item(4, 1)
Copy-Format $ws.Cells.Item(2, 1) $ws.Cells.Item(4, 2)
Copy-Format $ws.Cells.Item(2, 3) $ws.Cells.Item(4, 3)
Copy-Format $ws.Cells.Item(1, 4) $ws.Cells.Item(4, 4)
Copy-Format $ws.Cells.Item(2, 5) $ws.Cells.Item(4, 5)

# ---------------------------------------------------------------------------
# Row 5 — new row (A5/B5 stay empty but keep the header-style formatting)
# ---------------------------------------------------------------------------
$ws.Cells.Item(5, 3).Value = 100
$ws.Cells.Item(5, 4).Value = "'2022-07-22"
$ws.Cells.Item(5, 5).Value = "11:36"
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = 80
$ws.Cells.Item(5, 8).Value = 25

Copy-Format $ws.Cells.Item(2, 1) $ws.Cells.Item(5, 1)
Copy-Format $ws.Cells.Item(2, 2) $ws.Cells.Item(5, 2)
Copy-Format $ws.Cells.Item(2, 3) $ws.Cells.Item(5, 3)
Copy-Format $ws.Cells.Item(1, 4) $ws.Cells.Item(5, 4)
Copy-Format $ws.Cells.Item(2, 5) $ws.Cells.Item(5, 5)

# ---------------------------------------------------------------------------
# Merge the repeated name / wavelength cells down their groups
# ---------------------------------------------------------------------------
$ws.Range("A2:A5").Merge()
$ws.Range("B2:B3").Merge()
$ws.Range("B4:B5").Merge()
